# Update source plate concentrations/volumes (DA edit)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ECHO_source_plate")

# E2: the "Volume" note for the Water row changed first value 50 -> 19.475
$ws.Range("E2").Value = "19.475,50,50,65,65,65,65,65,65,65,65,65,65,65"

# Updated Concentration values (column E) for several wells
$ws.Range("E12").Value = 51.274999999999999
$ws.Range("E14").Value = 43.9
$ws.Range("E18").Value = 49.95
$ws.Range("E31").Value = 49.3
$ws.Range("E32").Value = 53.725000000000001
$ws.Range("E33").Value = 50.274999999999999
$ws.Range("E34").Value = 50.45
$ws.Range("E35").Value = 51.15

# Restore view state changes
$ws.Range("D23").Select()
